# Generate Report for Handoff
# Updates the localization-status report to reflect that the handoff
# package has been generated (status flips from "In Translation" to
# "Ready for handoff"), refreshes the two timestamp columns, and widens
# the Status/zh-cn/de-de columns so the longer status text fits.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn (E) and de-de (F) status + generate-date columns ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 18:38:36"

# --- zh-cn sheet: Status (C) + Latest Handoff Datetime (H) ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 18:38:31"

# --- de-de sheet: Status (C) + Latest Handoff Datetime (H) ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-28 18:38:36"

# --- Widen columns to fit the new, longer status text ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25
$wsZhCn.Columns.Item(3).ColumnWidth = 16.25
$wsDeDe.Columns.Item(3).ColumnWidth = 16.25
